# Apply scheduled-runner price/profit updates to the Leve profit sheets.
# Each sheet holds per-Leve pricing pulled from the market board; this
# batch refresh updates currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) for the rows whose underlying market data changed.
$wb = $excel.ActiveWorkbook


# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376
# Row 64
$ws.Range("H64").Value = 2870.2144
$ws.Range("I64").Value = 2535.7144
$ws.Range("J64").Value = 3204.7144
$ws.Range("K64").Value = 2535.7144
$ws.Range("L64").Value = 3204.7144
$ws.Range("M64").Value = -2287.7144
$ws.Range("N64").Value = -3700.7144
# Row 65
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880
# Row 67
$ws.Range("H67").Value = 2870.2144
$ws.Range("I67").Value = 2535.7144
$ws.Range("J67").Value = 3204.7144
$ws.Range("K67").Value = 2535.7144
$ws.Range("L67").Value = 3204.7144
$ws.Range("M67").Value = -1677.7144
$ws.Range("N67").Value = -4920.7144
# Row 70
$ws.Range("H70").Value = 1412.1428
$ws.Range("J70").Value = 1517
$ws.Range("L70").Value = 4551
$ws.Range("N70").Value = -5091
# Row 73
$ws.Range("H73").Value = 1412.1428
$ws.Range("J73").Value = 1517
$ws.Range("L73").Value = 4551
$ws.Range("N73").Value = -6423
# Row 76
$ws.Range("H76").Value = 1957151.4
$ws.Range("I76").Value = 2134202
$ws.Range("K76").Value = 2134202
$ws.Range("M76").Value = -2133887
# Row 79
$ws.Range("H79").Value = 1957151.4
$ws.Range("I79").Value = 2134202
$ws.Range("K79").Value = 2134202
$ws.Range("M79").Value = -2133110
# Row 86
$ws.Range("H86").Value = 5500
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
# Row 89
$ws.Range("H89").Value = 5500
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
# Row 96
$ws.Range("H96").Value = 540
$ws.Range("I96").Value = 540
$ws.Range("K96").Value = 1620
$ws.Range("M96").Value = -247
# Row 112
$ws.Range("H112").Value = 3788.8635
$ws.Range("J112").Value = 3874.0476
$ws.Range("L112").Value = 11622.1428
$ws.Range("N112").Value = -13838.1428
# Row 129
$ws.Range("H129").Value = 874.57776
$ws.Range("J129").Value = 873.9
$ws.Range("L129").Value = 2621.7
$ws.Range("N129").Value = -12621.7
# Row 138
$ws.Range("H138").Value = 3373.4546
$ws.Range("I138").Value = 6898.375
$ws.Range("K138").Value = 20695.125
$ws.Range("M138").Value = -15555.125

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2576.1729
$ws.Range("I32").Value = 1721.6232
$ws.Range("J32").Value = 7489.8335
$ws.Range("K32").Value = 1721.6232
$ws.Range("L32").Value = 7489.8335
$ws.Range("M32").Value = -1434.6232
$ws.Range("N32").Value = -8063.8335
# Row 45
$ws.Range("H45").Value = 1832.5555
# Row 97
$ws.Range("H97").Value = 937.3333
$ws.Range("I97").Value = 824.6667
$ws.Range("K97").Value = 824.6667
$ws.Range("M97").Value = -328.6667

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 655.2857
$ws.Range("I64").Value = 716
$ws.Range("J64").Value = 503.5
$ws.Range("K64").Value = 716
$ws.Range("L64").Value = 503.5
$ws.Range("M64").Value = -491
$ws.Range("N64").Value = -953.5
# Row 67
$ws.Range("H67").Value = 655.2857
$ws.Range("I67").Value = 716
$ws.Range("J67").Value = 503.5
$ws.Range("K67").Value = 716
$ws.Range("L67").Value = 503.5
$ws.Range("M67").Value = 64
$ws.Range("N67").Value = -2063.5
# Row 80
$ws.Range("H80").Value = 11584.777
$ws.Range("I80").Value = 498.5
$ws.Range("J80").Value = 14752.286
$ws.Range("K80").Value = 498.5
$ws.Range("L80").Value = 14752.286
$ws.Range("M80").Value = 499.5
$ws.Range("N80").Value = -16748.286
# Row 83
$ws.Range("H83").Value = 11584.777
$ws.Range("I83").Value = 498.5
$ws.Range("J83").Value = 14752.286
$ws.Range("K83").Value = 2492.5
$ws.Range("L83").Value = 73761.42999999999
$ws.Range("M83").Value = 2499.5
$ws.Range("N83").Value = -83745.42999999999
# Row 86
$ws.Range("H86").Value = 184710.45
$ws.Range("I86").Value = 2952.75
$ws.Range("K86").Value = 2952.75
$ws.Range("M86").Value = -1829.75
# Row 89
$ws.Range("H89").Value = 184710.45
$ws.Range("I89").Value = 2952.75
$ws.Range("K89").Value = 14763.75
$ws.Range("M89").Value = -9147.75
# Row 99
$ws.Range("H99").Value = 1554.5454
$ws.Range("I99").Value = 1262.5
$ws.Range("K99").Value = 1262.5
$ws.Range("M99").Value = 235.5

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 20
$ws.Range("H20").Value = 49998.332
$ws.Range("J20").Value = 49998.332
$ws.Range("L20").Value = 49998.332
$ws.Range("N20").Value = -50470.332
# Row 30
$ws.Range("H30").Value = 49998.332
$ws.Range("J30").Value = 49998.332
$ws.Range("L30").Value = 49998.332
$ws.Range("N30").Value = -50180.332
# Row 58
$ws.Range("H58").Value = 1448.7838
$ws.Range("I58").Value = 845.6539
$ws.Range("K58").Value = 845.6539
$ws.Range("M58").Value = -642.6539
# Row 128
$ws.Range("H128").Value = 49998.332
$ws.Range("J128").Value = 49998.332
$ws.Range("L128").Value = 49998.332
$ws.Range("N128").Value = -59958.332
# Row 132
$ws.Range("H132").Value = 2254.2778
$ws.Range("I132").Value = 1511.5238
$ws.Range("K132").Value = 4534.5714
$ws.Range("M132").Value = -2004.5714
# Row 134
$ws.Range("H134").Value = 1652.4
$ws.Range("I134").Value = 1417.258
$ws.Range("K134").Value = 4251.774
$ws.Range("M134").Value = -1716.774
# Row 136
$ws.Range("H136").Value = 1448.7838
$ws.Range("I136").Value = 845.6539
$ws.Range("K136").Value = 2536.9617
$ws.Range("M136").Value = 13.03830000000016

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 848.3889
$ws.Range("J107").Value = 945.6667
$ws.Range("L107").Value = 2837.0001
$ws.Range("N107").Value = -6677.0001
# Row 131
$ws.Range("H131").Value = 743.33673
$ws.Range("J131").Value = 757.837
$ws.Range("L131").Value = 2273.511
$ws.Range("N131").Value = -12353.511
# Row 139
$ws.Range("H139").Value = 499
$ws.Range("I139").Value = 499.5
$ws.Range("J139").Value = 498
$ws.Range("K139").Value = 1498.5
$ws.Range("L139").Value = 1494
$ws.Range("M139").Value = 3641.5
$ws.Range("N139").Value = -11774
# Row 140
$ws.Range("H140").Value = 1640.6842
$ws.Range("I140").Value = 925.1177
$ws.Range("J140").Value = 2219.9524
$ws.Range("K140").Value = 2775.3531
$ws.Range("L140").Value = 6659.8572
$ws.Range("M140").Value = 2404.6469
$ws.Range("N140").Value = -17019.8572

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 6897208.5
$ws.Range("I11").Value = 7113744.5
$ws.Range("K11").Value = 7113744.5
$ws.Range("M11").Value = -7113605.5
# Row 70
$ws.Range("H70").Value = 9534
$ws.Range("I70").Value = 11938.111
$ws.Range("K70").Value = 11938.111
$ws.Range("M70").Value = -11668.111
# Row 73
$ws.Range("H73").Value = 9534
$ws.Range("I73").Value = 11938.111
$ws.Range("K73").Value = 11938.111
$ws.Range("M73").Value = -11002.111
# Row 102
$ws.Range("H102").Value = 2863.9092
$ws.Range("I102").Value = 2850.3
$ws.Range("K102").Value = 2850.3
$ws.Range("M102").Value = -1228.3
# Row 107
$ws.Range("H107").Value = 1417.1666
# Row 132
$ws.Range("H132").Value = 2370.6216
$ws.Range("I132").Value = 1990.6072
$ws.Range("J132").Value = 3552.889
$ws.Range("K132").Value = 5971.821599999999
$ws.Range("L132").Value = 10658.667
$ws.Range("M132").Value = -3441.821599999999
$ws.Range("N132").Value = -15718.667

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3010.6875
$ws.Range("I7").Value = 2159.2727
$ws.Range("J7").Value = 4883.8
$ws.Range("K7").Value = 2159.2727
$ws.Range("L7").Value = 4883.8
$ws.Range("M7").Value = -2047.2727
$ws.Range("N7").Value = -5107.8
# Row 20
$ws.Range("H20").Value = 8671.5
$ws.Range("I20").Value = 6543
$ws.Range("J20").Value = 10800
$ws.Range("K20").Value = 6543
$ws.Range("L20").Value = 10800
$ws.Range("M20").Value = -6317
$ws.Range("N20").Value = -11252
# Row 22
$ws.Range("H22").Value = 4475
$ws.Range("I22").Value = 6200
$ws.Range("J22").Value = 2750
$ws.Range("K22").Value = 6200
$ws.Range("L22").Value = 2750
$ws.Range("M22").Value = -5905
$ws.Range("N22").Value = -3340
# Row 27
$ws.Range("H27").Value = 4475
$ws.Range("I27").Value = 6200
$ws.Range("J27").Value = 2750
$ws.Range("K27").Value = 6200
$ws.Range("L27").Value = 2750
$ws.Range("M27").Value = -6093
$ws.Range("N27").Value = -2964
# Row 126
$ws.Range("H126").Value = 3010.6875
$ws.Range("I126").Value = 2159.2727
$ws.Range("J126").Value = 4883.8
$ws.Range("K126").Value = 6477.8181
$ws.Range("L126").Value = 14651.4
$ws.Range("M126").Value = -4007.8181
$ws.Range("N126").Value = -19591.4
# Row 132
$ws.Range("H132").Value = 2466.0605
$ws.Range("I132").Value = 2155.923
$ws.Range("J132").Value = 2667.65
$ws.Range("K132").Value = 6467.768999999999
$ws.Range("L132").Value = 8002.950000000001
$ws.Range("M132").Value = -3937.768999999999
$ws.Range("N132").Value = -13062.95

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 108
$ws.Range("H108").Value = 64999
$ws.Range("J108").Value = 64999
$ws.Range("L108").Value = 64999
$ws.Range("N108").Value = -72679
# Row 130
$ws.Range("H130").Value = 34568.855
$ws.Range("J130").Value = 34568.855
$ws.Range("L130").Value = 34568.855
$ws.Range("N130").Value = -44608.855
# Row 132
$ws.Range("H132").Value = 1319.3158
$ws.Range("I132").Value = 1065.24
$ws.Range("J132").Value = 1807.9231
$ws.Range("K132").Value = 3195.72
$ws.Range("L132").Value = 5423.7693
$ws.Range("M132").Value = -665.7200000000003
$ws.Range("N132").Value = -10483.7693
